$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 3.5
$ws.Range("I2").Value = 2.45
$ws.Range("L2").Value = 3.5
$ws.Range("AB2").Value = 1.37
$ws.Range("AI2").Value = 4.5

# Row 3
$ws.Range("G3").Value = 3.1
$ws.Range("I3").Value = 2.1
$ws.Range("J3").Value = 3.6
$ws.Range("K3").Value = 2.25
$ws.Range("AA3").Value = 1.54
$ws.Range("AE3").Value = 12
$ws.Range("AH3").Value = 29
$ws.Range("AJ3").Value = 6.5
$ws.Range("AM3").Value = 9.5
$ws.Range("AS3").Value = 151

# Row 5
$ws.Range("G5").Value = 2.6
$ws.Range("H5").Value = 3.3
$ws.Range("I5").Value = 2.47
$ws.Range("J5").Value = 3.15
$ws.Range("X5").Value = 1.39
$ws.Range("AC5").Value = 10.25
$ws.Range("AD5").Value = 14.5
$ws.Range("AE5").Value = 9.75
$ws.Range("AF5").Value = 30
$ws.Range("AI5").Value = 11.75
$ws.Range("AJ5").Value = 6.6
$ws.Range("AO5").Value = 9.25
